$wb = $excel.ActiveWorkbook

# Rename sheets: Weekly -> Monthly
$wsGeneral = $wb.Worksheets.Item("GeneralTaxRateWeekly")
$wsGeneral.Name = "GeneralTaxRateMonthly"

$wsProcess = $wb.Worksheets.Item("ProcessPayrollForWeeklyTax")
$wsProcess.Name = "ProcessPayrollForMonthlyTax"

# Update references to the renamed sheets on the "first" worksheet
$wsFirst = $wb.Worksheets.Item("first")
$wsFirst.Range("A3").Value = "GeneralTaxRateMonthly"
$wsFirst.Range("A4").Value = "ProcessPayrollForMonthlyTax"

# Update employee marker text on each data worksheet (EMP 107 -> EMP 105)
$wsGeneral.Range("A2").Value = "DO NOT TOUCH AUTOMATION EMP 105"
$wsProcess.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 105"

$wsReports = $wb.Worksheets.Item("TestReports")
$wsReports.Range("B2").Value = "DO NOT TOUCH AUTOMATION EMP 105"

# Restore per-sheet selections, finishing with the sheet that should be active
$wsFirst.Range("F5").Select() | Out-Null
$wsProcess.Activate() | Out-Null
$wsProcess.Range("G9").Select() | Out-Null
$wsReports.Activate() | Out-Null
$wsReports.Range("J8").Select() | Out-Null
$wsGeneral.Activate() | Out-Null
$wsGeneral.Range("F16").Select() | Out-Null
